# Generate Report for Handoff
#
# Rows 4,5,7,8,9,10 (files 061e3db9, 3bc58189, 805e8a2f, aa57083f, d2921622,
# f2bb2f91) just had their handoff xliff regenerated:
#   - Overview!G  (Latest HO Xliff Generate Date)   -> 2016-11-14 07:39:06
#   - zh-cn!E     (Priority)                        -> ht
#   - zh-cn!H     (Latest Handoff Datetime)         -> 2016-11-14 07:38:51
#   - de-de!E     (Priority)                        -> ht
#   - de-de!H     (Latest Handoff Datetime)         -> 2016-11-14 07:39:06

$wb = $excel.ActiveWorkbook

$rows = @(4, 5, 7, 8, 9, 10)

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-11-14 07:39:06"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-11-14 07:38:51"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-11-14 07:39:06"
}
